$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
# This text appears on the "Overview" sheet (columns for zh-cn / de-de status)
# as well as on each locale sheet's own "Status" column.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        # NOTE: keep the literal on the left of -eq so PowerShell performs a
        # string comparison instead of coercing the literal to the type of
        # $v (some cells hold the text "True"/"False" which round-trips as
        # a System.Boolean, and "$true -eq <any non-empty string>" is $true).
        if ("Ready for handoff" -eq $v) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Shrink the "Status" columns to reflect the new, shorter text ---
# Overview sheet: columns E (zh-cn status) and F (de-de status)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
